$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update F2/G2 (Images / Description)
$ws.Range("F2").Value = "cintamani275_7"
$ws.Range("G2").Value = "Beschreibung"

# Row 3 - Buddhas - Manjusri
$ws.Range("B3").Value = "Buddhas - Manjusri"
$ws.Range("C3").Value = "manjusri"
$ws.Range("D3").Value = 649
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = "cintamani275_7"
$ws.Range("G3").Value = "Beschreibung"

# Row 4 - Malas - klein
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Malas  - klein"
$ws.Range("C4").Value = "kleine mala"
$ws.Range("D4").Value = 250
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = "cintamani275_7"
$ws.Range("G4").Value = "Beschreibung"

# Row 5 - Thangkas - groß
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Thangkas  - groß"
$ws.Range("C5").Value = "große Thangka"
$ws.Range("D5").Value = 250
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "cintamani275_7"
$ws.Range("G5").Value = "Beschreibung"

# Row 6 - Thangkas - groß 2
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Thangkas  - groß"
$ws.Range("C6").Value = "große Thangka 2"
$ws.Range("D6").Value = 250
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = "cintamani275_7"
$ws.Range("G6").Value = "Beschreibung"

# Row 7 - Thangkas - klein
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Thangkas  - klein"
$ws.Range("C7").Value = "kleine thangka"
$ws.Range("D7").Value = 250
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "cintamani275_7"
$ws.Range("G7").Value = "Beschreibung"

# Update selection to match target (F2 active cell)
[void]$ws.Range("F2").Select()
